$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the row containing ID=4 ("Matrimonio"), which is row 5.
# This removes the "Matrimonio" shared string, shifts the ID values of
# the rows below up into its place, and shifts every following row up
# by one.
$ws.Rows.Item(5).Delete()

# Leave the last-used cell selection further down the sheet, matching
# where the cursor ended up after the edit.
[void]$ws.Range("B15").Select()
